$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "81.414.25"
$ws.Range("E2").Value = "  +4.32%  "

# Row 3
$ws.Range("D3").Value = "3.179.67"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.95"
$ws.Range("E5").Value = "  +2.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.70"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.296"
$ws.Range("E7").Value = "  +29.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +3.08%  "

# Row 10
$ws.Range("D10").Value = "3.175.71"
$ws.Range("E10").Value = "  +0.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.590"
$ws.Range("E11").Value = "  +3.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000263"
$ws.Range("E12").Value = "  +16.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.166"
$ws.Range("E13").Value = "  +2.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.38"
$ws.Range("E14").Value = "  -1.12%  "

# Row 15
$ws.Range("D15").Value = "3.762.05"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.96"
$ws.Range("E16").Value = "  +1.05%  "

# Row 17
$ws.Range("D17").Value = "81.484.86"
$ws.Range("E17").Value = "  +4.56%  "

# Row 18
$ws.Range("D18").Value = "3.181.67"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.22"
$ws.Range("E19").Value = "  +13.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.19"
$ws.Range("E20").Value = "  -0.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.24"
$ws.Range("E21").Value = "  -1.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.79"
$ws.Range("E22").Value = "  +2.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("E23").Value = "  +6.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.16"
$ws.Range("E24").Value = "  +6.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.22"
$ws.Range("E25").Value = "  +10.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.24"
$ws.Range("E26").Value = "  +2.25%  "

# Row 27
$ws.Range("D27").Value = "3.343.94"
$ws.Range("E27").Value = "  +0.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.92"
$ws.Range("E28").Value = "  +0.72%  "

# Row 29
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000127"
$ws.Range("E30").Value = "  +10.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.16"
$ws.Range("E31").Value = "  +3.29%  "

# Row 32
$ws.Range("E32").Value = "  +0.62%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "561.05"
$ws.Range("E33").Value = "  +7.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  +2.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.04"
$ws.Range("E35").Value = "  +3.17%  "

# Row 36
$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.142"
$ws.Range("E36").Value = "  +31.45%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  +12.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.13"
$ws.Range("E38").Value = "  +2.89%  "

# Row 39
$ws.Range("E39").Value = "  +0.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.415"
$ws.Range("E40").Value = "  +4.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.12"
$ws.Range("E41").Value = "  +21.94%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.97"
$ws.Range("E42").Value = "  +10.59%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.04"
$ws.Range("E43").Value = "  +15.37%  "

# Row 44
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.76"
$ws.Range("E44").Value = "  +3.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "160.34"
$ws.Range("E45").Value = "  -2.07%  "

# Row 46
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "189.87"
$ws.Range("E47").Value = "  -3.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.35"
$ws.Range("E48").Value = "  +4.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.50"
$ws.Range("E49").Value = "  +3.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.785"
$ws.Range("E50").Value = "  -1.63%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.28"
$ws.Range("E51").Value = "  +4.78%  "
